$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for rows 2-8 from
# 45207 (2023-10-08) to 45208 (2023-10-09).
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 45208
}
